# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for the
# leves whose market data refreshed in this scheduled run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 259
$ws.Range("I5").Value = 65.5
$ws.Range("J5").Value = 336.4
$ws.Range("K5").Value = 65.5
$ws.Range("L5").Value = 336.4
$ws.Range("M5").Value = 49.5
$ws.Range("N5").Value = -566.4

# Row 42
$ws.Range("H42").Value = 5830
$ws.Range("I42").Value = 1000
$ws.Range("J42").Value = 6635
$ws.Range("K42").Value = 3000
$ws.Range("L42").Value = 19905
$ws.Range("M42").Value = -2770
$ws.Range("N42").Value = -20365

# Row 101
$ws.Range("H101").Value = 25004718
$ws.Range("J101").Value = 1000
$ws.Range("L101").Value = 3000
$ws.Range("N101").Value = -6244

# Row 112
$ws.Range("H112").Value = 4739.8887
$ws.Range("J112").Value = 4882.375
$ws.Range("L112").Value = 14647.125
$ws.Range("N112").Value = -16863.125

# Row 123
$ws.Range("H123").Value = 28999
$ws.Range("J123").Value = 28999
$ws.Range("L123").Value = 28999
$ws.Range("N123").Value = -38799

$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1074
$ws.Range("I94").Value = 1074
$ws.Range("K94").Value = 1074
$ws.Range("M94").Value = -623

# Row 102
$ws.Range("H102").Value = 65618.664
$ws.Range("I102").Value = 65618.664
$ws.Range("K102").Value = 65618.664
$ws.Range("M102").Value = -62373.664

# Row 107
$ws.Range("H107").Value = 1852.6
$ws.Range("I107").Value = 1852.6
$ws.Range("K107").Value = 1852.6
$ws.Range("M107").Value = 67.40000000000009

$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 49780
$ws.Range("J20").Value = 49780
$ws.Range("L20").Value = 49780
$ws.Range("N20").Value = -50252

# Row 30
$ws.Range("H30").Value = 49780
$ws.Range("J30").Value = 49780
$ws.Range("L30").Value = 49780
$ws.Range("N30").Value = -49962

# Row 68
$ws.Range("H68").Value = 99
$ws.Range("J68").Value = 99
$ws.Range("L68").Value = 99
$ws.Range("N68").Value = -1597

# Row 71
$ws.Range("H71").Value = 99
$ws.Range("J71").Value = 99
$ws.Range("L71").Value = 297
$ws.Range("N71").Value = -7785

# Row 112
$ws.Range("H112").Value = 40702
$ws.Range("J112").Value = 40702
$ws.Range("L112").Value = 40702
$ws.Range("N112").Value = -43656

# Row 128
$ws.Range("H128").Value = 49780
$ws.Range("J128").Value = 49780
$ws.Range("L128").Value = 49780
$ws.Range("N128").Value = -59740

$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 2186.111
$ws.Range("J22").Value = 2271.875
$ws.Range("L22").Value = 6815.625
$ws.Range("N22").Value = -7153.625

# Row 23
$ws.Range("H23").Value = 650.875
$ws.Range("I23").Value = 149
$ws.Range("J23").Value = 818.1667
$ws.Range("K23").Value = 447
$ws.Range("L23").Value = 2454.5001
$ws.Range("M23").Value = -212
$ws.Range("N23").Value = -2924.5001

# Row 27
$ws.Range("H27").Value = 2186.111
$ws.Range("J27").Value = 2271.875
$ws.Range("L27").Value = 6815.625
$ws.Range("N27").Value = -7019.625

# Row 64
$ws.Range("H64").Value = 1522.5
$ws.Range("I64").Value = 696.6667
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 2090.0001
$ws.Range("L64").Value = 12000
$ws.Range("M64").Value = -1820.0001
$ws.Range("N64").Value = -12540

# Row 67
$ws.Range("H67").Value = 1522.5
$ws.Range("I67").Value = 696.6667
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 2090.0001
$ws.Range("L67").Value = 12000
$ws.Range("M67").Value = -1154.0001
$ws.Range("N67").Value = -13872

# Row 75
$ws.Range("H75").Value = 3562.5
$ws.Range("J75").Value = 3083.3333
$ws.Range("L75").Value = 9249.999899999999
$ws.Range("N75").Value = -11245.9999

# Row 78
$ws.Range("H78").Value = 3562.5
$ws.Range("J78").Value = 3083.3333
$ws.Range("L78").Value = 27749.9997
$ws.Range("N78").Value = -37733.9997

# Row 86
$ws.Range("H86").Value = 594
$ws.Range("J86").Value = 1000
$ws.Range("L86").Value = 3000
$ws.Range("N86").Value = -5372

# Row 89
$ws.Range("H89").Value = 594
$ws.Range("J89").Value = 1000
$ws.Range("L89").Value = 9000
$ws.Range("N89").Value = -20856

# Row 98
$ws.Range("H98").Value = 1199.091
$ws.Range("J98").Value = 1463.4286
$ws.Range("L98").Value = 4390.2858
$ws.Range("N98").Value = -7386.2858

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 2672
$ws.Range("I70").Value = 2672
$ws.Range("K70").Value = 2672
$ws.Range("M70").Value = -2402

# Row 73
$ws.Range("H73").Value = 2672
$ws.Range("I73").Value = 2672
$ws.Range("K73").Value = 2672
$ws.Range("M73").Value = -1736

# Row 80
$ws.Range("H80").Value = 2701.2
$ws.Range("J80").Value = 2876.5
$ws.Range("L80").Value = 2876.5
$ws.Range("N80").Value = -4872.5

# Row 83
$ws.Range("H83").Value = 2701.2
$ws.Range("J83").Value = 2876.5
$ws.Range("L83").Value = 14382.5
$ws.Range("N83").Value = -24366.5

# Row 122
$ws.Range("H122").Value = 6755.1665
$ws.Range("I122").Value = 6106.6
$ws.Range("K122").Value = 18319.8
$ws.Range("M122").Value = -15869.8

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 9187.615
$ws.Range("I7").Value = 12565
$ws.Range("K7").Value = 12565
$ws.Range("M7").Value = -12453

# Row 46
$ws.Range("H46").Value = 2114.9
$ws.Range("I46").Value = 1174.5
$ws.Range("J46").Value = 2350
$ws.Range("K46").Value = 1174.5
$ws.Range("L46").Value = 2350
$ws.Range("M46").Value = -986.5
$ws.Range("N46").Value = -2726

# Row 68
$ws.Range("H68").Value = 1132.7333
$ws.Range("J68").Value = 1899.4
$ws.Range("L68").Value = 1899.4
$ws.Range("N68").Value = -3397.4

# Row 71
$ws.Range("H71").Value = 1132.7333
$ws.Range("J71").Value = 1899.4
$ws.Range("L71").Value = 9497
$ws.Range("N71").Value = -16985

# Row 100
$ws.Range("H100").Value = 2933.6875
$ws.Range("I100").Value = 2788.5715
$ws.Range("J100").Value = 3949.5
$ws.Range("K100").Value = 2788.5715
$ws.Range("L100").Value = 3949.5
$ws.Range("M100").Value = -2247.5715
$ws.Range("N100").Value = -5031.5

# Row 126
$ws.Range("H126").Value = 9187.615
$ws.Range("I126").Value = 12565
$ws.Range("K126").Value = 37695
$ws.Range("M126").Value = -35225

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 19357.428
$ws.Range("I62").Value = 19875.5
$ws.Range("J62").Value = 18666.666
$ws.Range("K62").Value = 19875.5
$ws.Range("L62").Value = 18666.666
$ws.Range("M62").Value = -19251.5
$ws.Range("N62").Value = -19914.666

# Row 65
$ws.Range("H65").Value = 19357.428
$ws.Range("I65").Value = 19875.5
$ws.Range("J65").Value = 18666.666
$ws.Range("K65").Value = 99377.5
$ws.Range("L65").Value = 93333.33
$ws.Range("M65").Value = -96257.5
$ws.Range("N65").Value = -99573.33

# Row 80
$ws.Range("H80").Value = 287
$ws.Range("I80").Value = 273
$ws.Range("J80").Value = 301
$ws.Range("K80").Value = 273
$ws.Range("L80").Value = 301
$ws.Range("M80").Value = 725
$ws.Range("N80").Value = -2297

# Row 83
$ws.Range("H83").Value = 287
$ws.Range("I83").Value = 273
$ws.Range("J83").Value = 301
$ws.Range("K83").Value = 819
$ws.Range("L83").Value = 903
$ws.Range("M83").Value = 4173
$ws.Range("N83").Value = -10887

# Row 132
$ws.Range("H132").Value = 1416.3334
$ws.Range("I132").Value = 1416.3334
$ws.Range("K132").Value = 4249.0002
$ws.Range("M132").Value = -1719.0002
